$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value assignments (matching the target diff)
$values = @{
    "D2" = 14002
    "E2" = 720
    "F2" = 720
    "G2" = 627
    "H2" = 506
    "I2" = 506
    "K2" = 16259
    "L2" = 12678
    "M2" = 3581
    "N2" = 3581
    "P2" = 1000
    "Q2" = 635
    "R2" = -749
    "S2" = -281
    "T2" = 547
    "U2" = 88
    "V2" = 1070
    "W2" = 5.14
    "X2" = 3.61
    "Y2" = 14.58
    "Z2" = 3.28
    "AA2" = 354.02
    "AB2" = 258.06
    "AC2" = 2528
    "AE2" = 17906
    "AF2" = 0
    "AG2" = 630
    "AI2" = 24.92
    "AJ2" = 20000000
    "D3" = 19037
    "E3" = 1122
    "F3" = 1122
    "G3" = 1049
    "H3" = 830
    "I3" = 830
    "K3" = 17543
    "L3" = 11772
    "M3" = 5771
    "N3" = 5771
    "P3" = 1100
    "Q3" = 90
    "R3" = -649
    "S3" = 866
    "T3" = 409
    "U3" = -319
    "V3" = 555
    "W3" = 5.89
    "X3" = 4.36
    "Y3" = 17.76
    "Z3" = 4.91
    "AA3" = 203.99
    "AB3" = 424.51
    "AC3" = 4043
    "AD3" = 25.6
    "AE3" = 26231
    "AF3" = 3.95
    "AG3" = 940
    "AH3" = 0.91
    "AI3" = 24.91
    "AJ3" = 22000000
    "D4" = 18608
    "E4" = 876
    "F4" = 876
    "G4" = 834
    "H4" = 771
    "I4" = 771
    "K4" = 17915
    "L4" = 11598
    "M4" = 6316
    "N4" = 6316
    "P4" = 1100
    "Q4" = -1077
    "R4" = -996
    "S4" = 1844
    "T4" = 822
    "U4" = -1898
    "V4" = 2606
    "W4" = 4.71
    "X4" = 4.14
    "Y4" = 12.75
    "Z4" = 4.35
    "AA4" = 183.62
    "AB4" = 474.16
    "AC4" = 3502
    "AD4" = 22.98
    "AE4" = 28711
    "AF4" = 2.8
    "AG4" = 940
    "AH4" = 1.17
    "AI4" = 26.84
    "AJ4" = 22000000
    "D5" = 17613
    "E5" = 43
    "F5" = 43
    "G5" = -251
    "H5" = -86
    "I5" = -86
    "K5" = 20068
    "L5" = 14028
    "M5" = 6040
    "N5" = 6040
    "P5" = 1100
    "Q5" = -364
    "R5" = -1150
    "S5" = 2285
    "T5" = 1029
    "U5" = -1393
    "V5" = 5099
    "W5" = 0.24
    "X5" = -0.49
    "Y5" = -1.4
    "Z5" = -0.45
    "AA5" = 232.26
    "AB5" = 449.02
    "AC5" = -392
    "AD5" = -152.48
    "AE5" = 27454
    "AF5" = 2.18
    "AG5" = 500
    "AH5" = 0.84
    "AI5" = -127.49
    "AJ5" = 22000000
    "AG6" = 500
    "AH6" = 1.35
    "AJ6" = 22000000
    "D7" = 15651
    "E7" = 342
    "G7" = 259
    "H7" = 180
    "I7" = 180
    "K7" = 24270
    "L7" = 17700
    "M7" = 6570
    "N7" = 6570
    "P7" = 1100
    "Q7" = 920
    "R7" = -290
    "S7" = -20
    "T7" = 590
    "U7" = 340
    "W7" = 2.19
    "X7" = 1.15
    "Y7" = 2.85
    "Z7" = 0.84
    "AA7" = 269.41
    "AC7" = 818
    "AD7" = 34.04
    "AE7" = 29864
    "AF7" = 0.93
    "AG7" = 511
    "AH7" = 1.84
    "AI7" = 62.47
    "D8" = 17272
    "E8" = 624
    "G8" = 498
    "H8" = 450
    "I8" = 373
    "K8" = 26330
    "L8" = 19420
    "M8" = 6910
    "N8" = 6910
    "P8" = 1100
    "Q8" = 1550
    "R8" = -620
    "S8" = 60
    "T8" = 600
    "U8" = 950
    "W8" = 3.61
    "X8" = 2.61
    "Y8" = 5.53
    "Z8" = 1.78
    "AA8" = 281.04
    "AC8" = 1695
    "AD8" = 16.43
    "AE8" = 31409
    "AF8" = 0.89
    "AG8" = 528
    "AH8" = 1.9
    "AI8" = 31.13
    "D9" = 20025
    "E9" = 1025
    "G9" = 900
    "H9" = 670
    "I9" = 735
    "K9" = 28990
    "L9" = 21520
    "M9" = 7470
    "N9" = 7470
    "P9" = 1100
    "Q9" = 1940
    "R9" = -650
    "S9" = 140
    "T9" = 620
    "U9" = 1320
    "W9" = 5.12
    "X9" = 3.35
    "Y9" = 10.22
    "Z9" = 2.42
    "AA9" = 288.09
    "AC9" = 3341
    "AD9" = 8.34
    "AE9" = 33955
    "AF9" = 0.82
    "AG9" = 550
    "AH9" = 1.97
    "AI9" = 16.46
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# Cells that are removed entirely in the target (no replacement value)
$clears = @(
    "J2", "O2", "AD2", "AH2", "J3", "O3", "J4", "O4", "J5", "O5", "D6", "E6", "F6", "G6", "H6", "I6", "K6", "L6", "M6", "N6", "P6", "Q6", "R6", "S6", "T6", "U6", "V6", "W6", "X6", "Y6", "Z6", "AA6", "AB6", "AC6", "AD6", "AE6", "AF6", "AI6"
)

foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

Write-Output "Applied $($values.Count) value updates and $($clears.Count) clears."